$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Scénarii - F2" (1st sheet) : shrink the "Interface" block down to a
# single new sub-scenario ("Interface IHM") and drop the two now-unused
# trailing scenario blocks.
# ---------------------------------------------------------------------------
$wsF2 = $wb.Worksheets.Item(1)

# Old row 10 (a plain data row, "3") is not needed any more: drop it. This
# shifts the old (empty) merged sub-header placeholder row up from 11 to 10,
# reusing its existing merge/style instead of creating a brand new one.
$wsF2.Rows(10).Delete()

# Drop the now-unused trailing rows (everything after the new row 11, which
# is the old row "12" - the first row of the next block - shifted up)
$wsF2.Rows("12:17").Delete()

$wsF2.Range("B10").Value = "Interface IHM"
$wsF2.Range("C11").Value = "Le manipulateur clique sur l'application Qt et décide de lancer ou non un essai"

$wsF2.Columns("C").ColumnWidth = 82.83072916666667
$wsF2.Range("F10").Select()

# ---------------------------------------------------------------------------
# Sheet "Scénarii - F5" (4th sheet) : fill in the missing step descriptions
# and drop the 3 trailing scenario blocks, keeping only a single, completed
# one.
# ---------------------------------------------------------------------------
$wsF5 = $wb.Worksheets.Item(4)

# Drop old rows 9-19 (three trailing scenario blocks)
$wsF5.Rows("9:19").Delete()

# Old row 8 was an (empty) merged sub-header placeholder - unmerge it so it
# can become the closing data row of the remaining block
$wsF5.Range("B8:E8").UnMerge()

$wsF5.Range("C5").Value = "Le manipulateur clique sur le bouton pour commencer un essai"
$wsF5.Range("C6").Value = "Les données sont visible par le manipulateur "
$wsF5.Range("C7").Value = "Le manipulateur clique sur le bouton ""courbe"" pour voir les points de la courbe"

# Row 8 becomes the last data row of the group, so it needs the "closing"
# (top+bottom only) border variant of the normal row style
$wsF5.Range("B5").Copy()
$wsF5.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsF5.Range("B8").Borders(10).LineStyle = -4142
$wsF5.Range("B8").Value = 4

$wsF5.Range("E5").Copy()
$wsF5.Range("E8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsF5.Range("E8").Borders(7).LineStyle = -4142

$wsF5.Range("C5").Copy()
$wsF5.Range("C8").PasteSpecial(-4122)
$wsF5.Range("D8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$wsF5.Range("C8").Borders(7).LineStyle = -4142
$wsF5.Range("C8").Borders(10).LineStyle = -4142
$wsF5.Range("D8").Borders(7).LineStyle = -4142
$wsF5.Range("D8").Borders(10).LineStyle = -4142
$wsF5.Range("C8").Value = "Le manipulateur visualise la même courbe mais de type oscilloscope"

$wsF5.Columns("C").ColumnWidth = 91.83072916666667
$wsF5.Range("C19").Select()

# "Scénarii - F5" becomes the active/visible tab
$wsF5.Activate()

$excel.CutCopyMode = 0
